# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Palta" (avocado) at the top of the
# "Feria Lagunitas de Puerto Montt" data block (rows 226-227), pushing the
# existing rows down by two. The new rows reuse the same market / category
# metadata as the block and only carry new date + price information.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 226 - existing row 226 (and everything below)
# shifts down to row 228, etc.
$ws.Rows.Item(226).Insert()
$ws.Rows.Item(226).Insert()

# --- New row 226 ---------------------------------------------------------
$ws.Cells.Item(226, 1).Value = 4
$ws.Cells.Item(226, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(226, 3).Value = "Los Lagos"
$ws.Cells.Item(226, 4).Value = "12/13/2021"
$ws.Cells.Item(226, 5).Value = 10
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100106
$ws.Cells.Item(226, 8).Value = "Oleaginosos"
$ws.Cells.Item(226, 9).Value = 100106002
$ws.Cells.Item(226, 10).Value = "Palta"
$ws.Cells.Item(226, 11).Value = "Hass"
$ws.Cells.Item(226, 12).Value = "Primera"
$ws.Cells.Item(226, 13).Value = 200
$ws.Cells.Item(226, 14).Value = 4000
$ws.Cells.Item(226, 15).Value = 4100
$ws.Cells.Item(226, 16).Value = 4050
$ws.Cells.Item(226, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(226, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(226, 19).Value = 4050
$ws.Cells.Item(226, 20).Value = 1

# --- New row 227 ---------------------------------------------------------
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = "12/13/2021"
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100106
$ws.Cells.Item(227, 8).Value = "Oleaginosos"
$ws.Cells.Item(227, 9).Value = 100106002
$ws.Cells.Item(227, 10).Value = "Palta"
$ws.Cells.Item(227, 11).Value = "Hass"
$ws.Cells.Item(227, 12).Value = "Segunda"
$ws.Cells.Item(227, 13).Value = 100
$ws.Cells.Item(227, 14).Value = 3600
$ws.Cells.Item(227, 15).Value = 3600
$ws.Cells.Item(227, 16).Value = 3600
$ws.Cells.Item(227, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(227, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(227, 19).Value = 3600
$ws.Cells.Item(227, 20).Value = 1
